# Auto-generated edit script applying Mateus_Profits.xlsx market-data refresh
# Updates computed market columns (H:N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4798.2856
$ws.Range("I40").Value = 2980.6667
$ws.Range("J40").Value = 7221.778
$ws.Range("K40").Value = 2980.6667
$ws.Range("L40").Value = 7221.778
$ws.Range("M40").Value = -2805.6667
$ws.Range("N40").Value = -7571.778

$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15540
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -16872
$ws.Range("M73").ClearContents()

$ws.Range("H86").Value = 3051.5334
$ws.Range("I86").Value = 2834
$ws.Range("J86").Value = 3649.75
$ws.Range("K86").Value = 2834
$ws.Range("L86").Value = 3649.75
$ws.Range("M86").Value = -1711
$ws.Range("N86").Value = -5895.75

$ws.Range("H89").Value = 3051.5334
$ws.Range("I89").Value = 2834
$ws.Range("J89").Value = 3649.75
$ws.Range("K89").Value = 14170
$ws.Range("L89").Value = 18248.75
$ws.Range("M89").Value = -8554
$ws.Range("N89").Value = -29480.75

$ws.Range("H100").Value = 1362.579
$ws.Range("I100").Value = 1029.2142
$ws.Range("K100").Value = 1029.2142
$ws.Range("M100").Value = -488.2141999999999

$ws.Range("H137").Value = 2343.1365
$ws.Range("I137").Value = 2160.0667
$ws.Range("J137").Value = 2735.4285
$ws.Range("K137").Value = 6480.2001
$ws.Range("L137").Value = 8206.2855
$ws.Range("M137").Value = -3930.2001
$ws.Range("N137").Value = -13306.2855

$ws.Range("H138").Value = 2796.3953
$ws.Range("I138").Value = 2398.5
$ws.Range("J138").Value = 2916.9697
$ws.Range("K138").Value = 7195.5
$ws.Range("L138").Value = 8750.909100000001
$ws.Range("M138").Value = -2055.5
$ws.Range("N138").Value = -19030.9091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2435.111
$ws.Range("I2").Value = 1783.2
$ws.Range("K2").Value = 1783.2
$ws.Range("M2").Value = -1670.2

$ws.Range("H32").Value = 9548
$ws.Range("I32").Value = 9548
$ws.Range("K32").Value = 9548
$ws.Range("M32").Value = -9261

$ws.Range("H63").Value = 1965.4166
$ws.Range("I63").Value = 1755
$ws.Range("J63").Value = 2596.6667
$ws.Range("K63").Value = 1755
$ws.Range("L63").Value = 2596.6667
$ws.Range("M63").Value = -1069
$ws.Range("N63").Value = -3968.6667

$ws.Range("H66").Value = 1965.4166
$ws.Range("I66").Value = 1755
$ws.Range("J66").Value = 2596.6667
$ws.Range("K66").Value = 8775
$ws.Range("L66").Value = 12983.3335
$ws.Range("M66").Value = -5343
$ws.Range("N66").Value = -19847.3335

$ws.Range("H74").Value = 5821.0435
$ws.Range("I74").Value = 4430.25
$ws.Range("K74").Value = 4430.25
$ws.Range("M74").Value = -3556.25

$ws.Range("H77").Value = 5821.0435
$ws.Range("I77").Value = 4430.25
$ws.Range("K77").Value = 22151.25
$ws.Range("M77").Value = -17783.25

$ws.Range("H97").Value = 967.6
$ws.Range("I97").Value = 999.2143
$ws.Range("K97").Value = 999.2143
$ws.Range("M97").Value = -503.2143

$ws.Range("H116").Value = 2435.111
$ws.Range("I116").Value = 1783.2
$ws.Range("K116").Value = 1783.2
$ws.Range("M116").Value = 510.8

$ws.Range("H132").Value = 2554.4048
$ws.Range("I132").Value = 2519.8647
$ws.Range("J132").Value = 2810
$ws.Range("K132").Value = 7559.5941
$ws.Range("L132").Value = 8430
$ws.Range("M132").Value = -5029.5941
$ws.Range("N132").Value = -13490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 69189.5
$ws.Range("J2").Value = 69189.5
$ws.Range("L2").Value = 69189.5
$ws.Range("N2").Value = -69415.5

$ws.Range("H3").Value = 2435.111
$ws.Range("I3").Value = 1783.2
$ws.Range("K3").Value = 1783.2
$ws.Range("M3").Value = -1669.2

$ws.Range("H105").Value = 573.5714
$ws.Range("I105").Value = 521.25
$ws.Range("K105").Value = 521.25
$ws.Range("M105").Value = 1225.75

$ws.Range("H134").Value = 3444.6584
$ws.Range("I134").Value = 3405.775
$ws.Range("K134").Value = 10217.325
$ws.Range("M134").Value = -7682.325000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4361.4614
$ws.Range("I16").Value = 2358.5715
$ws.Range("K16").Value = 2358.5715
$ws.Range("M16").Value = -2071.5715

$ws.Range("H31").Value = 5806.95
$ws.Range("J31").Value = 5756.357
$ws.Range("L31").Value = 5756.357
$ws.Range("N31").Value = -6346.357

$ws.Range("H34").Value = 5806.95
$ws.Range("J34").Value = 5756.357
$ws.Range("L34").Value = 5756.357
$ws.Range("N34").Value = -6160.357

$ws.Range("H58").Value = 4713.5405
$ws.Range("I58").Value = 3076.08
$ws.Range("K58").Value = 3076.08
$ws.Range("M58").Value = -2873.08

$ws.Range("H74").Value = 42432
$ws.Range("J74").Value = 42432
$ws.Range("L74").Value = 42432
$ws.Range("N74").Value = -44180

$ws.Range("H77").Value = 42432
$ws.Range("J77").Value = 42432
$ws.Range("L77").Value = 127296
$ws.Range("N77").Value = -136032

$ws.Range("H86").Value = 53990.363
$ws.Range("I86").Value = 8748
$ws.Range("K86").Value = 8748
$ws.Range("M86").Value = -7625

$ws.Range("H89").Value = 53990.363
$ws.Range("I89").Value = 8748
$ws.Range("K89").Value = 43740
$ws.Range("M89").Value = -38124

$ws.Range("H99").Value = 4397.636
$ws.Range("I99").Value = 4711
$ws.Range("J99").Value = 3849.25
$ws.Range("K99").Value = 4711
$ws.Range("L99").Value = 3849.25
$ws.Range("M99").Value = -3213
$ws.Range("N99").Value = -6845.25

$ws.Range("H113").Value = 4361.4614
$ws.Range("I113").Value = 2358.5715
$ws.Range("K113").Value = 2358.5715
$ws.Range("M113").Value = -188.5715

$ws.Range("H126").Value = 4397.636
$ws.Range("I126").Value = 4711
$ws.Range("J126").Value = 3849.25
$ws.Range("K126").Value = 14133
$ws.Range("L126").Value = 11547.75
$ws.Range("M126").Value = -11663
$ws.Range("N126").Value = -16487.75

$ws.Range("H136").Value = 4713.5405
$ws.Range("I136").Value = 3076.08
$ws.Range("K136").Value = 9228.24
$ws.Range("M136").Value = -6678.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 75.71429000000001
$ws.Range("I2").Value = 72.5
$ws.Range("K2").Value = 435
$ws.Range("M2").Value = -322

$ws.Range("H4").Value = 7135464
$ws.Range("I4").Value = 13816392
$ws.Range("K4").Value = 41449176
$ws.Range("M4").Value = -41449064

$ws.Range("H28").Value = 4115
$ws.Range("J28").Value = 3700
$ws.Range("L28").Value = 11100
$ws.Range("N28").Value = -11564

$ws.Range("H37").Value = 125273.5
$ws.Range("J37").Value = 125273.5
$ws.Range("L37").Value = 375820.5
$ws.Range("N37").Value = -376044.5

$ws.Range("H38").Value = 1036.1111
$ws.Range("I38").Value = 30
$ws.Range("J38").Value = 1161.875
$ws.Range("K38").Value = 90
$ws.Range("L38").Value = 3485.625
$ws.Range("M38").Value = 257
$ws.Range("N38").Value = -4179.625

$ws.Range("H138").Value = 1763
$ws.Range("I138").Value = 1763
$ws.Range("K138").Value = 5289
$ws.Range("M138").Value = -149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 38894
$ws.Range("J93").Value = 38894
$ws.Range("L93").Value = 38894
$ws.Range("N93").Value = -42638

$ws.Range("H99").Value = 7785.5
$ws.Range("I99").Value = 6342.6
$ws.Range("K99").Value = 6342.6
$ws.Range("M99").Value = -4096.6

$ws.Range("H102").Value = 3079.875
$ws.Range("I102").Value = 1645.2941
$ws.Range("K102").Value = 1645.2941
$ws.Range("M102").Value = -23.29410000000007

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1287.2609
$ws.Range("I16").Value = 1365.0667
$ws.Range("K16").Value = 1365.0667
$ws.Range("M16").Value = -1195.0667

$ws.Range("H22").Value = 1598.2858
$ws.Range("I22").Value = 1559.6
$ws.Range("J22").Value = 1695
$ws.Range("K22").Value = 1559.6
$ws.Range("L22").Value = 1695
$ws.Range("M22").Value = -1264.6
$ws.Range("N22").Value = -2285

$ws.Range("H27").Value = 1598.2858
$ws.Range("I27").Value = 1559.6
$ws.Range("J27").Value = 1695
$ws.Range("K27").Value = 1559.6
$ws.Range("L27").Value = 1695
$ws.Range("M27").Value = -1452.6
$ws.Range("N27").Value = -1909

$ws.Range("H61").Value = 101799.3
$ws.Range("I61").Value = 126873.5
$ws.Range("K61").Value = 126873.5
$ws.Range("M61").Value = -126671.5

$ws.Range("H113").Value = 101799.3
$ws.Range("I113").Value = 126873.5
$ws.Range("K113").Value = 126873.5
$ws.Range("M113").Value = -124703.5

$ws.Range("H132").Value = 15168.393
$ws.Range("I132").Value = 19964.375
$ws.Range("J132").Value = 8773.75
$ws.Range("K132").Value = 59893.125
$ws.Range("L132").Value = 26321.25
$ws.Range("M132").Value = -57363.125
$ws.Range("N132").Value = -31381.25

$ws.Range("H136").Value = 4684.875
$ws.Range("I136").Value = 4709.2
$ws.Range("J136").Value = 4598
$ws.Range("K136").Value = 14127.6
$ws.Range("L136").Value = 13794
$ws.Range("M136").Value = -11577.6
$ws.Range("N136").Value = -18894

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 47160
$ws.Range("J103").Value = 47160
$ws.Range("L103").Value = 47160
$ws.Range("N103").Value = -49504
